$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the two cells that held "first/待定" to the new text "first/时间待定"
$ws.Range("B11").Value = "first/时间待定"
$ws.Range("B12").Value = "first/时间待定"
